$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new BOM row 33: Neudynium Magnets
$ws.Range("A33").Value = "Neudynium Magnets"
$ws.Range("B33").Value = 2
$ws.Range("C33").Value = "Any"
$ws.Range("D33").Value = "5 mm OD x 2 mm W"

# Match the wrap-text style already used throughout column D
$ws.Range("D33").WrapText = $true

# Update the saved selection / view state to match the edited workbook
$ws.Range("D28").Select()
